$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2022 column (S) to the table, mirroring the existing column R
# formatting (borders, fonts, number format) by copying its formats over.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 76.099999999999994

$excel.CutCopyMode = $false

# Move the selection to P8 to match the new view state.
$ws.Range("P8").Select()
